$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2018 LEAVE CREDITS")

# Row 42: period end date shifts from 1/1/2023 to 1/31/2023
$ws.Range("A42").Value = 44957

# Row 43: becomes a formula (EDATE one month after A42) instead of a literal date
$ws.Range("A43").Formula = "=EDATE(A42,1)"

# Row 44: period end date shifts from 3/1/2023 to 3/31/2023
$ws.Range("A44").Value = 45016

# Row 46: period end date shifts from 4/1/2023 to 4/30/2023, and earns 1.25 days
$ws.Range("A46").Value = 45046
$ws.Range("C46").Value = 1.25

# Row 47: new period end date 5/31/2023, earns 1.25 days
$ws.Range("A47").Value = 45077
$ws.Range("C47").Value = 1.25

# Row 48: new period end date 6/30/2023, remarks "SP(1-0-0)", and a date in K48
$ws.Range("A48").Value = 45107
$ws.Range("B48").Value = "SP(1-0-0)"
$ws.Range("K48").Value = 45112
$ws.Range("K45").Copy()
$ws.Range("K48").PasteSpecial(-4122)

# Rows 49-56: fill in period end dates for subsequent months
$ws.Range("A49").Value = 45138
$ws.Range("A50").Value = 45169
$ws.Range("A51").Value = 45199
$ws.Range("A52").Value = 45230
$ws.Range("A53").Value = 45260
$ws.Range("A54").Value = 45291
$ws.Range("A55").Value = 45322
$ws.Range("A56").Value = 45351
